# Re-solve the first-order ODE dy/dt = y on the interval [0, 2]
# (was [0, 4]) with the same number of sample points (20), using
# Euler's method, mirroring `numpy.linspace` + a forward-Euler loop,
# and write the refreshed t/y table back into A2:B21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$t0 = 0
$tf = 2
$n = 20
$div = $n - 1
$step = ($tf - $t0) / $div

# Build the t samples exactly like numpy.linspace: t[i] = i*step + t0,
# with the final sample snapped to the exact endpoint.
$t = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $n; $i++) {
    $t[$i] = $i * $step + $t0
}
$t[$n - 1] = $tf

$h = $t[1] - $t[0]

$y = 1
for ($i = 0; $i -lt $n; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $t[$i]
    $ws.Cells.Item($row, 2).Value = $y
    $y = $y + $h * $y
}

Write-Output "updated $n rows"
